$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is digits-and-dots text that Excel would
# otherwise auto-convert to a number; force Text format, assign,
# then restore the cell style so no lasting formatting change remains.
$textForceCells = @("D5", "D6", "D10", "D12", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D48")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.009.99"
$ws.Range("E2").Value = "  +2.57%  "
$ws.Range("D3").Value = "2.959.22"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "573.33"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "160.23"
$ws.Range("E6").Value = "  +5.17%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").Value = "2.958.25"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "6.62"
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "66.146.44"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").Value = "3.453.55"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "2.959.53"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "445.63"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "13.67"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "0.672"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").Value = "7.11"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "81.82"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "12.14"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "9.95"
$ws.Range("E28").Value = "  -9.44%  "
$ws.Range("D29").Value = "8.05"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("E30").Value = "  +7.78%  "
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").Value = "0.0₃0980"
$ws.Range("E32").Value = "  -10.02%  "
$ws.Range("D33").Value = "27.03"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "0.977"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "5.67"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "49.24"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "1.98"
$ws.Range("E39").Value = "  -5.57%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "0.298"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("B41").Value = "Arweave"
$ws.Range("C41").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D41").Value = "43.02"
$ws.Range("E41").Value = "  -2.18%  "
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "2.80"
$ws.Range("E43").Value = "  -8.75%  "
$ws.Range("D44").Value = "8.32"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").Value = "379.86"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "2.704.21"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").Value = "131.19"
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("E51").Value = "  +4.21%  "

foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
